$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = 1.509
$ws.Range("C2").Value = "BALLENOIL"
$ws.Range("D2").Value = "AVENIDA FUENLABRADA, 6"
$ws.Range("E2").Value = "HUMANES DE MADRID"
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = 1.519
$ws.Range("C3").Value = "SHELL"
$ws.Range("D3").Value = "AVENIDA FUENLABRADA, 110"
$ws.Range("E3").Value = "HUMANES DE MADRID"
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = 1.639
$ws.Range("C4").Value = "BALLENOIL"
$ws.Range("D4").Value = "CALLE CABO RUFINO LAZARO, 7"
$ws.Range("E4").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = 1.679
$ws.Range("C5").Value = "GALP"
$ws.Range("D5").Value = "CTRA. N-VI km 21,700"
$ws.Range("E5").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = 1.699
$ws.Range("C6").Value = "CARREFOUR"
$ws.Range("D6").Value = "CARRETERA MADRID-LA CORUÑA KM. 22"
$ws.Range("E6").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = 1.709
$ws.Range("C7").Value = "REPSOL"
$ws.Range("D7").Value = "CL MADRID, 52"
$ws.Range("E7").Value = "HUMANES DE MADRID"
$ws.Range("A8").Value = "7"
$ws.Range("B8").Value = 1.718
$ws.Range("C8").Value = "REPSOL"
$ws.Range("D8").Value = "CALLE COPENHAGUES/N, S/N"
$ws.Range("E8").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A9").Value = "8"
$ws.Range("B9").Value = 1.718
$ws.Range("C9").Value = "REPSOL"
$ws.Range("D9").Value = "A-6 km 25,5"
$ws.Range("E9").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A10").Value = "9"
$ws.Range("B10").Value = 1.718
$ws.Range("C10").Value = "REPSOL"
$ws.Range("D10").Value = "CTRA. M-505 km 5,5"
$ws.Range("E10").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A11").Value = "10"
$ws.Range("B11").Value = 1.718
$ws.Range("C11").Value = "REPSOL"
$ws.Range("D11").Value = "CARRETERA M-505 km 5.5"
$ws.Range("E11").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A12").Value = "11"
$ws.Range("B12").Value = 1.719
$ws.Range("C12").Value = "BP VALDONAIRE"
$ws.Range("D12").Value = "CARRETERA AVD.DE LA INDUSTRIA KM. 15"
$ws.Range("E12").Value = "HUMANES DE MADRID"
$ws.Range("A13").Value = "12"
$ws.Range("B13").Value = 1.719
$ws.Range("C13").Value = "BP HUMANES - EL MOLINO"
$ws.Range("D13").Value = "AVENIDA DE LAS FLORES, 2"
$ws.Range("E13").Value = "HUMANES DE MADRID"
$ws.Range("A14").Value = "13"
$ws.Range("B14").Value = 1.719
$ws.Range("C14").Value = "CEPSA"
$ws.Range("D14").Value = "CARRETERA M-405 KM. 5,6"
$ws.Range("E14").Value = "HUMANES DE MADRID"
$ws.Range("A15").Value = "14"
$ws.Range("B15").Value = 1.719
$ws.Range("C15").Value = "REPSOL HUMANES"
$ws.Range("D15").Value = "AVENIDA LA INDUSTRIA, S/N"
$ws.Range("E15").Value = "HUMANES DE MADRID"
$ws.Range("A16").Value = "15"
$ws.Range("B16").Value = 1.719
$ws.Range("C16").Value = "REPSOL"
$ws.Range("D16").Value = "CARRETERA AVENIDA  DE LA INDUSTRIA , 46 KM. 1,1"
$ws.Range("E16").Value = "HUMANES DE MADRID"
$ws.Range("A17").Value = "16"
$ws.Range("B17").Value = 1.719
$ws.Range("C17").Value = "BP LAS ROZAS"
$ws.Range("D17").Value = "CL LAS CRUCES  S/N"
$ws.Range("E17").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A18").Value = "17"
$ws.Range("B18").Value = 1.719
$ws.Range("C18").Value = "REPSOL"
$ws.Range("D18").Value = "CR A-6, 20,3"
$ws.Range("E18").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A19").Value = "18"
$ws.Range("B19").Value = 1.725
$ws.Range("C19").Value = "CEPSA"
$ws.Range("D19").Value = "CARRETERA M-405 KM. 6"
$ws.Range("E19").Value = "HUMANES DE MADRID"
$ws.Range("A20").Value = "19"
$ws.Range("B20").Value = 1.924
$ws.Range("C20").Value = "COSTCO"
$ws.Range("D20").Value = "CALLE INNOVACIÓN, 19"
$ws.Range("E20").Value = "ROZAS DE MADRID (LAS)"
$ws.Range("A21").Value = "20"
$ws.Range("C21").Value = "T9"
$ws.Range("D21").Value = "CALLE TENERIFE (DE), 2"
$ws.Range("E21").Value = "HUMANES DE MADRID"